$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the ESD protection block (Zener diodes D3-D10 -> D3-D12) ---
# Reference designator list grows from 8 diodes to 10 (D11, D12 added).
$ws.Range("B21").Value = "D3, D4, D5, D6, D7, D8, D9, D10, D11, D12"

# Quantity used goes from 8 to 10 to match the added reference designators.
# (Cost per-row and the grand total below recalculate automatically.)
$ws.Range("H21").Value = 10

# --- Normalize leftover one-off cell formatting back to the standard style ---
# B3, B4 and B13 previously had a stray "wrap text" style; B21 previously had a
# stray highlight fill. Copy plain formatting from a normally-styled cell
# (B5 already uses the standard, unhighlighted, non-wrapped style) onto them,
# without touching their existing text content.
$ws.Range("B5").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
